# Highlights the "dubious / uncertain" requirement blocks of the
# Specifiche Supplementari document, per the commit message:
#   "Specifiche supplementari in cui sono evidenziate le cose dubbi o incerti"
#
# For each subsection heading below, the heading paragraph and the body
# paragraph right after it get marked with a highlight colour:
#   - "Fattori umani" (under "Usabilità")      -> yellow
#   - "Performance"                             -> green
#   - "Configurabilità"                         -> yellow
#
# WdColorIndex values: wdYellow = 7, wdBrightGreen = 4
$wdYellow = 7
$wdBrightGreen = 4

$d = $word.ActiveDocument

function Get-ParagraphIndexByText($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $text) {
            return $i
        }
    }
    return -1
}

function Set-HeadingHighlight($paraIndex, $colorIndex) {
    # Heading paragraphs already carry paragraph-mark run properties
    # (w:pPr/w:rPr, e.g. bold). Going through the Font object propagates
    # the highlight both to the run(s) AND to those paragraph-mark run
    # properties - matching Word's behaviour when the whole paragraph
    # (pilcrow included) is selected and highlighted.
    $d.Paragraphs($paraIndex).Range.Font.HighlightColorIndex = $colorIndex
}

function Set-BodyHighlight($paraIndex, $colorIndex) {
    # Body paragraphs have no paragraph-mark run properties of their own;
    # only their run(s) should pick up the highlight, so set the property
    # straight on the Range (not via Font) to avoid minting a new w:pPr.
    $d.Paragraphs($paraIndex).Range.HighlightColorIndex = $colorIndex
}

function Set-SectionHighlight($headingText, $colorIndex) {
    $headingIdx = Get-ParagraphIndexByText($headingText)
    Set-HeadingHighlight $headingIdx $colorIndex
    Set-BodyHighlight ($headingIdx + 1) $colorIndex
}

Set-SectionHighlight "Usabilità" $wdYellow
Set-SectionHighlight "Fattori umani" $wdYellow
Set-SectionHighlight "Performance" $wdBrightGreen
Set-SectionHighlight "Configurabilità" $wdYellow
